# Auto-generated edit script
# Applies updated currentAveragePrice / Leve profit values across the Exodus_Profits workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 941.03076
$ws.Range("I15").Value = 941.03076
$ws.Range("K15").Value = 2823.09228
$ws.Range("M15").Value = -2654.09228
$ws.Range("H19").Value = 508.76923
$ws.Range("J19").Value = 536.3333
$ws.Range("L19").Value = 536.3333
$ws.Range("N19").Value = -886.3333
$ws.Range("H74").Value = 3593.5715
$ws.Range("I74").Value = 3150.5881
$ws.Range("J74").Value = 5476.25
$ws.Range("K74").Value = 3150.5881
$ws.Range("L74").Value = 5476.25
$ws.Range("M74").Value = -2214.5881
$ws.Range("N74").Value = -7348.25
$ws.Range("H77").Value = 3593.5715
$ws.Range("I77").Value = 3150.5881
$ws.Range("J77").Value = 5476.25
$ws.Range("K77").Value = 15752.9405
$ws.Range("L77").Value = 27381.25
$ws.Range("M77").Value = -11072.9405
$ws.Range("N77").Value = -36741.25
$ws.Range("H98").Value = 1299.7
$ws.Range("I98").Value = 1299.7
$ws.Range("K98").Value = 1299.7
$ws.Range("M98").Value = 198.3
$ws.Range("H122").Value = 1299.7
$ws.Range("I122").Value = 1299.7
$ws.Range("K122").Value = 3899.1
$ws.Range("M122").Value = -1449.1
$ws.Range("H131").Value = 1866.0
$ws.Range("I131").Value = 799.0
$ws.Range("K131").Value = 2397.0
$ws.Range("M131").Value = 2643.0
$ws.Range("H132").Value = 1879.1212
$ws.Range("I132").Value = 1440.069
$ws.Range("K132").Value = 4320.207
$ws.Range("M132").Value = -1790.207
$ws.Range("H137").Value = 764616.94
$ws.Range("I137").Value = 1028.0
$ws.Range("K137").Value = 3084.0
$ws.Range("M137").Value = -534.0
$ws.Range("H138").Value = 10640082.0
$ws.Range("I138").Value = 1569.4166
$ws.Range("J138").Value = 12196937.0
$ws.Range("K138").Value = 4708.2498
$ws.Range("L138").Value = 36590811.0
$ws.Range("M138").Value = 431.7502000000004
$ws.Range("N138").Value = -36601091.0

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7677.8853
$ws.Range("I32").Value = 4544.6387
$ws.Range("K32").Value = 4544.6387
$ws.Range("M32").Value = -4257.6387
$ws.Range("H63").Value = 11614.667
$ws.Range("I63").Value = 2422.5
$ws.Range("K63").Value = 2422.5
$ws.Range("M63").Value = -1736.5
$ws.Range("H66").Value = 11614.667
$ws.Range("I66").Value = 2422.5
$ws.Range("K66").Value = 12112.5
$ws.Range("M66").Value = -8680.5
$ws.Range("H88").Value = 715658.56
$ws.Range("I88").Value = 1251197.0
$ws.Range("K88").Value = 1251197.0
$ws.Range("M88").Value = -1250791.0
$ws.Range("H91").Value = 715658.56
$ws.Range("I91").Value = 1251197.0
$ws.Range("K91").Value = 1251197.0
$ws.Range("M91").Value = -1249793.0
$ws.Range("H122").Value = 3939.2
$ws.Range("I122").Value = 3935.5454
$ws.Range("K122").Value = 11806.6362
$ws.Range("M122").Value = -9356.6362
$ws.Range("H132").Value = 2231.2444
$ws.Range("I132").Value = 1921.4706
$ws.Range("K132").Value = 5764.4118
$ws.Range("M132").Value = -3234.4118

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2491.5625
$ws.Range("I86").Value = 1677.2
$ws.Range("K86").Value = 1677.2
$ws.Range("M86").Value = -554.2
$ws.Range("H89").Value = 2491.5625
$ws.Range("I89").Value = 1677.2
$ws.Range("K89").Value = 8386.0
$ws.Range("M89").Value = -2770.0
$ws.Range("H94").Value = 1392.0769
$ws.Range("I94").Value = 1091.4166
$ws.Range("K94").Value = 1091.4166
$ws.Range("M94").Value = -640.4166
$ws.Range("H134").Value = 1970.6061
$ws.Range("I134").Value = 1575.2593
$ws.Range("K134").Value = 4725.7779
$ws.Range("M134").Value = -2190.7779

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3109.7307
$ws.Range("I31").Value = 2229.353
$ws.Range("J31").Value = 4772.6665
$ws.Range("K31").Value = 2229.353
$ws.Range("L31").Value = 4772.6665
$ws.Range("M31").Value = -1934.353
$ws.Range("N31").Value = -5362.6665
$ws.Range("H34").Value = 3109.7307
$ws.Range("I34").Value = 2229.353
$ws.Range("J34").Value = 4772.6665
$ws.Range("K34").Value = 2229.353
$ws.Range("L34").Value = 4772.6665
$ws.Range("M34").Value = -2027.353
$ws.Range("N34").Value = -5176.6665
$ws.Range("H58").Value = 2140.5
$ws.Range("I58").Value = 1895.6154
$ws.Range("J58").Value = 2352.7334
$ws.Range("K58").Value = 1895.6154
$ws.Range("L58").Value = 2352.7334
$ws.Range("M58").Value = -1692.6154
$ws.Range("N58").Value = -2758.7334
$ws.Range("H122").Value = 3474.625
$ws.Range("I122").Value = 2810.4443
$ws.Range("K122").Value = 8431.332900000001
$ws.Range("M122").Value = -5981.332900000001
$ws.Range("H132").Value = 2022.2174
$ws.Range("I132").Value = 957.8889
$ws.Range("K132").Value = 2873.6667
$ws.Range("M132").Value = -343.6667000000002
$ws.Range("H134").Value = 1345.95
$ws.Range("I134").Value = 826.89746
$ws.Range("J134").Value = 2309.9048
$ws.Range("K134").Value = 2480.69238
$ws.Range("L134").Value = 6929.714399999999
$ws.Range("M134").Value = 54.30762000000004
$ws.Range("N134").Value = -11999.7144
$ws.Range("H136").Value = 2140.5
$ws.Range("I136").Value = 1895.6154
$ws.Range("J136").Value = 2352.7334
$ws.Range("K136").Value = 5686.8462
$ws.Range("L136").Value = 7058.2002
$ws.Range("M136").Value = -3136.8462
$ws.Range("N136").Value = -12158.2002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 160.0
$ws.Range("I14").Value = 160.0
$ws.Range("K14").Value = 480.0
$ws.Range("M14").Value = -307.0
$ws.Range("H54").Value = 3398.3333
$ws.Range("J54").Value = 4995.0
$ws.Range("L54").Value = 14985.0
$ws.Range("N54").Value = -16103.0
$ws.Range("H64").Value = 2625.0
$ws.Range("I64").Value = 0.0
$ws.Range("J64").Value = 2625.0
$ws.Range("K64").Value = 0.0
$ws.Range("L64").Value = 7875.0
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -8415.0
$ws.Range("H67").Value = 2625.0
$ws.Range("I67").Value = 0.0
$ws.Range("J67").Value = 2625.0
$ws.Range("K67").Value = 0.0
$ws.Range("L67").Value = 7875.0
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -9747.0
$ws.Range("H75").Value = 16933.666
$ws.Range("I75").Value = 400.0
$ws.Range("K75").Value = 1200.0
$ws.Range("M75").Value = -202.0
$ws.Range("H78").Value = 16933.666
$ws.Range("I78").Value = 400.0
$ws.Range("K78").Value = 3600.0
$ws.Range("M78").Value = 1392.0
$ws.Range("H126").Value = 5005.077
$ws.Range("I126").Value = 4507.5
$ws.Range("J126").Value = 6663.6665
$ws.Range("K126").Value = 13522.5
$ws.Range("L126").Value = 19990.9995
$ws.Range("M126").Value = -8582.5
$ws.Range("N126").Value = -29870.9995
$ws.Range("H131").Value = 424812.78
$ws.Range("I131").Value = 91861.63
$ws.Range("J131").Value = 668976.94
$ws.Range("K131").Value = 275584.89
$ws.Range("L131").Value = 2006930.82
$ws.Range("M131").Value = -270544.89
$ws.Range("N131").Value = -2017010.82

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3298.7778
$ws.Range("I80").Value = 3338.8
$ws.Range("J80").Value = 3248.75
$ws.Range("K80").Value = 3338.8
$ws.Range("L80").Value = 3248.75
$ws.Range("M80").Value = -2340.8
$ws.Range("N80").Value = -5244.75
$ws.Range("H83").Value = 3298.7778
$ws.Range("I83").Value = 3338.8
$ws.Range("J83").Value = 3248.75
$ws.Range("K83").Value = 16694.0
$ws.Range("L83").Value = 16243.75
$ws.Range("M83").Value = -11702.0
$ws.Range("N83").Value = -26227.75
$ws.Range("H119").Value = 69449.336
$ws.Range("J119").Value = 69449.336
$ws.Range("L119").Value = 69449.336
$ws.Range("N119").Value = -79125.336
$ws.Range("H122").Value = 6541.433
$ws.Range("I122").Value = 6278.577
$ws.Range("K122").Value = 18835.731
$ws.Range("M122").Value = -16385.731
$ws.Range("H132").Value = 1712.2632
$ws.Range("I132").Value = 1298.5555
$ws.Range("K132").Value = 3895.6665
$ws.Range("M132").Value = -1365.6665
$ws.Range("H140").Value = 66708.57
$ws.Range("J140").Value = 79990.0
$ws.Range("L140").Value = 79990.0
$ws.Range("N140").Value = -90350.0

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1738165.4
$ws.Range("I40").Value = 1773.76
$ws.Range("K40").Value = 1773.76
$ws.Range("M40").Value = -1637.76
$ws.Range("H101").Value = 17440.5
$ws.Range("J101").Value = 17440.5
$ws.Range("L101").Value = 17440.5
$ws.Range("N101").Value = -23930.5
$ws.Range("H122").Value = 8719974.0
$ws.Range("I122").Value = 34987.0
$ws.Range("K122").Value = 104961.0
$ws.Range("M122").Value = -102511.0
$ws.Range("H132").Value = 2757.4583
$ws.Range("I132").Value = 2058.5715
$ws.Range("K132").Value = 6175.7145
$ws.Range("M132").Value = -3645.7145

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H94").Value = 37098.9
$ws.Range("I94").Value = 33331.668
$ws.Range("J94").Value = 38713.43
$ws.Range("K94").Value = 33331.668
$ws.Range("L94").Value = 38713.43
$ws.Range("M94").Value = -32430.668
$ws.Range("N94").Value = -40515.43
$ws.Range("H103").Value = 25000.0
$ws.Range("J103").Value = 25000.0
$ws.Range("L103").Value = 25000.0
$ws.Range("N103").Value = -27344.0
$ws.Range("H132").Value = 1210406.4
$ws.Range("I132").Value = 2313.8
$ws.Range("K132").Value = 6941.400000000001
$ws.Range("M132").Value = -4411.400000000001
$ws.Range("H136").Value = 3211.1428
$ws.Range("I136").Value = 2132.9062
$ws.Range("K136").Value = 6398.7186
$ws.Range("M136").Value = -3848.7186

Write-Host "Updated 243 cells and cleared 2 cells across 8 sheets."
